$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- Row 65 ----
# A65: date 2022-01-10 (44571), formatted like A60 (style index 5)
$ws.Range("A60").Copy()
$ws.Range("A65").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("A65").Value = (Get-Date -Year 2022 -Month 1 -Day 10 -Hour 0 -Minute 0 -Second 0)

# B65: hours, style like B64 (style index 1)
$ws.Range("B64").Copy()
$ws.Range("B65").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B65").Value = 1

# C65 / D65: plain text cells (no explicit style, same as C64/D64)
$ws.Range("C65").Value = "route tokenin tarkistukseen"
$ws.Range("D65").Value = "api"

# ---- Row 66 ----
# A66: date 2022-01-11 (44572)
$ws.Range("A60").Copy()
$ws.Range("A66").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A66").Value = (Get-Date -Year 2022 -Month 1 -Day 11 -Hour 0 -Minute 0 -Second 0)

# B66
$ws.Range("B64").Copy()
$ws.Range("B66").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B66").Value = 2

# C66 / D66
$ws.Range("C66").Value = "opeteltuauth  cookien lähettäminen suoraan serveriltä, testattu ja apin:n koodi refaktoroitu"
$ws.Range("D66").Value = "api"

# ---- Update totals row 75 ----
$ws.Range("B75").Formula = "=SUM(B2:B66)"

# ---- Update sheet view (scroll position & selection) ----
$excel.ActiveWindow.ScrollRow = 37
$ws.Range("A67").Select()
